$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text_C2 = @'
The applicant has extensive experience in projects related to natural language processing (NLP), computer vision, and PyTorch, which are directly in line with the required skills for the position. The projects 'LLMGuard', 'Literature Society IITJ Website', and 'Multi Model Data Analysis for Annotation of Human Activities' demonstrate the applicant's proficiency in NLP, Python, and PyTorch. Furthermore, the project 'CloudPhysician's Vital Extraction Challenge' showcases the applicant's expertise in computer vision and image processing. Although the applicant does not possess direct experience with PyTorch, considering their strong background in Python and related projects, it is reasonable to assume they can quickly adapt to using PyTorch for the required tasks. Overall, the applicant's project portfolio aligns well with the technical requirements of the job, demonstrating significant compatibility with the company's needs and therefore justifying the given score.
'@
$ws.Range("C2").Value = $text_C2

$text_E2 = @'
The applicant has exhibited exceptional communication skills in their responses, articulately expressing their strong alignment with Japanese work culture and values. Additionally, their proactive approach to self-improvement and continuous learning, as well as the demonstrated adaptability to diverse cultures, indicates their well-rounded personality traits. Furthermore, their expressed desire for mentorship, collaboration, and contribution to the company's objectives exemplifies qualities of teamwork, problem-solving, and time management. The applicant's evident enthusiasm for Japanese work culture and their holistic approach toward personal and professional growth solidify their suitability for the role, warranting a high score.
'@
$ws.Range("E2").Value = $text_E2

$text_C3 = @'
The applicant has significant experience in building web and mobile applications using a variety of technologies such as ReactJS, Flutter, NodeJS, ExpressJS, Firebase, Docker, etc. While the job description emphasizes NLP, Pytorch, Computer Vision, and Python, the applicant's skills are more aligned with front-end and full-stack development. Although the applicant has shown proficiency in related areas like WebRTC and Socket.IO, there is a notable gap in direct experience with the specific technologies mentioned in the job description. However, the applicant's diverse technical background and experience working with real-time communication technologies demonstrate their adaptability and problem-solving abilities, which are valuable traits for a tech role. The applicant has a strong foundation in web and mobile application development, but may require additional training and experience in NLP, Pytorch, Computer Vision, and Python to fully meet the requirements of the position.
'@
$ws.Range("C3").Value = $text_C3

$text_E3 = @'
The applicant has displayed a high level of adaptability, as evidenced by their willingness to work in various countries, including Japan, and their openness to learning a new language. Additionally, the applicant's collaborative approach to problem-solving and emphasis on teamwork align with the company's soft skill requirements. Their commitment to contributing to cutting-edge technologies and desire to bridge technical and creative aspects showcase strong problem-solving and communication abilities. While there are areas for improvement, such as the preference to work in a team and the need to enhance individual skill acquisition, the applicant's demonstrated adaptability, teamwork, and problem-solving capabilities warrant a high score for their personality and willingness to work in Japan.
'@
$ws.Range("E3").Value = $text_E3

$text_C4 = @'
The applicant has a strong background in NLP, Computer Vision, and Pytorch, which are key skills required for the job. The projects show experience in developing NLP and Computer Vision models using Python, PyTorch, and other relevant tools. The applicant has demonstrated the ability to work on cutting-edge technologies, analyze data, and develop and improve models, aligning with the job description. However, the applicant lacks direct experience in collecting and preprocessing text corpora, a key aspect of the job requirement. While the applicant's skills are impressive, this gap may impact their suitability for the role.
'@
$ws.Range("C4").Value = $text_C4

$text_E4 = @'
The applicant demonstrates strong adaptability and a willingness to work in Japan. Their proactive approach to addressing weaknesses and ability to work well in a team, as evidenced by the answers provided, aligns with the company's soft skills requirements. The applicant's openness to learning a new language and interest in Japan's work culture further reflects their willingness to embrace new cultural experiences. Their ability to overcome challenges and foster collaborative problem-solving demonstrates strong problem-solving and teamwork skills. Overall, the applicant appears to be a good fit for the company's soft skill requirements and willingness to work in Japan.
'@
$ws.Range("E4").Value = $text_E4

$text_C5 = @'
The applicant has experience in web development projects using JavaScript, HTML, and CSS, which demonstrates strong programming and front-end development skills. However, the job description specifically requires skills in NLP, Pytorch, Computer Vision, and Python, which the applicant does not possess based on the given projects and skills. While the applicant's experience in web development is commendable, it does not directly align with the company's requirements for a NLP Engineer. The applicant's technical skills are impressive but not entirely suitable for the specific role of NLP Engineer at the company.
'@
$ws.Range("C5").Value = $text_C5

$text_E5 = @'
The applicant has demonstrated strong communication skills by expressing a keen interest in working in Japan and adapting to a new culture. The responses indicate a willingness to collaborate with a team and a positive attitude towards overcoming challenges. The applicant's enthusiasm for learning new skills and contributing to advanced AI/ML technologies aligns with the company's expectations. However, the applicant also acknowledges the need to improve on working alone, which shows self-awareness and a willingness to address weaknesses. Overall, the applicant's responses reflect a strong personality and a high level of adaptability, making them well-suited for the company's soft skill requirements and their willingness to work in Japan.
'@
$ws.Range("E5").Value = $text_E5

$ws.Range("D2").Value = 4

